$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2").Value = 0.6
